$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Summary  (aggregate stats bumped because 3 more leadlag trades closed
# and a new leadlag trade was opened)
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 44
$summary.Range("E2").Value = "+10.8736%"
$summary.Range("F2").Value = "+0.2471%"

$summary.Range("C3").Value = 51
$summary.Range("D3").Value = "39.2%"
$summary.Range("E3").Value = "+6.6720%"
$summary.Range("F3").Value = "+0.1308%"

# ---------------------------------------------------------------------------
# Sheet: leadlag  (rows 32-34 transition from OPEN to CLOSED via the 5 minute
# time exit, and a brand new trade #64 is appended as row 53)
# ---------------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

# Row 32 - trade #42
$leadlag.Cells.Item(32, 7).Value = 69713.270189
$leadlag.Cells.Item(32, 8).Value = "CLOSED"
$leadlag.Cells.Item(32, 9).Value = 1.3956
$leadlag.Cells.Item(32, 10).Value = 13.96
$leadlag.Cells.Item(32, 13).Value = "time_exit_5min"
$leadlag.Cells.Item(32, 14).Value = 5

# Row 33 - trade #43
$leadlag.Cells.Item(33, 7).Value = 69546.00107499999
$leadlag.Cells.Item(33, 8).Value = "CLOSED"
$leadlag.Cells.Item(33, 9).Value = 1.1732
$leadlag.Cells.Item(33, 10).Value = 11.73
$leadlag.Cells.Item(33, 13).Value = "time_exit_5min"
$leadlag.Cells.Item(33, 14).Value = 5

# Row 34 - trade #44
$leadlag.Cells.Item(34, 7).Value = 69133.467752
$leadlag.Cells.Item(34, 8).Value = "CLOSED"
$leadlag.Cells.Item(34, 9).Value = -0.6047
$leadlag.Cells.Item(34, 10).Value = -6.05
$leadlag.Cells.Item(34, 13).Value = "time_exit_5min"
$leadlag.Cells.Item(34, 14).Value = 5

# Row 53 (new) - trade #64, freshly opened
$leadlag.Cells.Item(53, 1).Value = 64
$leadlag.Cells.Item(53, 2).NumberFormat = "@"
$leadlag.Cells.Item(53, 2).Value = "2026-02-16"
$leadlag.Cells.Item(53, 3).NumberFormat = "@"
$leadlag.Cells.Item(53, 3).Value = "21:34:53"
$leadlag.Cells.Item(53, 4).Value = "leadlag"
$leadlag.Cells.Item(53, 5).Value = "DOWN"
$leadlag.Cells.Item(53, 6).Value = 68688.94
$leadlag.Cells.Item(53, 8).Value = "OPEN"
$leadlag.Cells.Item(53, 9).Value = 0
$leadlag.Cells.Item(53, 10).Value = 0
$leadlag.Cells.Item(53, 11).Value = 0.75
$leadlag.Cells.Item(53, 12).Value = "Binance leading with -0.082% move"
$leadlag.Cells.Item(53, 14).Value = 0

# ---------------------------------------------------------------------------
# Sheet: All Trades  (mirrors the three newly closed leadlag trades)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Row 43 - trade #42
$allTrades.Cells.Item(43, 1).Value = 42
$allTrades.Cells.Item(43, 2).NumberFormat = "@"
$allTrades.Cells.Item(43, 2).Value = "2026-02-16"
$allTrades.Cells.Item(43, 3).NumberFormat = "@"
$allTrades.Cells.Item(43, 3).Value = "21:29:24"
$allTrades.Cells.Item(43, 4).Value = "leadlag"
$allTrades.Cells.Item(43, 5).Value = "UP"
$allTrades.Cells.Item(43, 6).Value = 68753.72
$allTrades.Cells.Item(43, 7).Value = 69713.270189
$allTrades.Cells.Item(43, 8).Value = "CLOSED"
$allTrades.Cells.Item(43, 9).Value = 1.3956
$allTrades.Cells.Item(43, 10).Value = 13.96
$allTrades.Cells.Item(43, 11).Value = 0.75
$allTrades.Cells.Item(43, 12).Value = "Binance leading with 0.157% move"
$allTrades.Cells.Item(43, 13).Value = "time_exit_5min"
$allTrades.Cells.Item(43, 14).Value = 5

# Row 44 - trade #43
$allTrades.Cells.Item(44, 1).Value = 43
$allTrades.Cells.Item(44, 2).NumberFormat = "@"
$allTrades.Cells.Item(44, 2).Value = "2026-02-16"
$allTrades.Cells.Item(44, 3).NumberFormat = "@"
$allTrades.Cells.Item(44, 3).Value = "21:29:31"
$allTrades.Cells.Item(44, 4).Value = "leadlag"
$allTrades.Cells.Item(44, 5).Value = "UP"
$allTrades.Cells.Item(44, 6).Value = 68739.52499999999
$allTrades.Cells.Item(44, 7).Value = 69546.00107499999
$allTrades.Cells.Item(44, 8).Value = "CLOSED"
$allTrades.Cells.Item(44, 9).Value = 1.1732
$allTrades.Cells.Item(44, 10).Value = 11.73
$allTrades.Cells.Item(44, 11).Value = 0.677
$allTrades.Cells.Item(44, 12).Value = "Coinbase leading with 0.068% move"
$allTrades.Cells.Item(44, 13).Value = "time_exit_5min"
$allTrades.Cells.Item(44, 14).Value = 5

# Row 45 - trade #44
$allTrades.Cells.Item(45, 1).Value = 44
$allTrades.Cells.Item(45, 2).NumberFormat = "@"
$allTrades.Cells.Item(45, 2).Value = "2026-02-16"
$allTrades.Cells.Item(45, 3).NumberFormat = "@"
$allTrades.Cells.Item(45, 3).Value = "21:29:37"
$allTrades.Cells.Item(45, 4).Value = "leadlag"
$allTrades.Cells.Item(45, 5).Value = "DOWN"
$allTrades.Cells.Item(45, 6).Value = 68717.925
$allTrades.Cells.Item(45, 7).Value = 69133.467752
$allTrades.Cells.Item(45, 8).Value = "CLOSED"
$allTrades.Cells.Item(45, 9).Value = -0.6047
$allTrades.Cells.Item(45, 10).Value = -6.05
$allTrades.Cells.Item(45, 11).Value = 0.604
$allTrades.Cells.Item(45, 12).Value = "Binance leading with -0.060% move"
$allTrades.Cells.Item(45, 13).Value = "time_exit_5min"
$allTrades.Cells.Item(45, 14).Value = 5

# ---------------------------------------------------------------------------
# Sheet: Comparison  (leadlag aggregate row)
# ---------------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Range("B2").Value = 51
$comparison.Range("C2").Value = "39.2%"
$comparison.Range("D2").Value = "2.49"
$comparison.Range("E2").Value = "+0.5580%"
$comparison.Range("F2").Value = "-0.3452%"
$comparison.Range("G2").Value = "1.62"
$comparison.Range("H2").Value = "-0.6047%"

$wb.Save()
